# Apply the updated crypto price/volume snapshot (GitHub Actions refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Cells.Item(2, 4).Value = "42.968.49"
$ws.Cells.Item(2, 5).Value = "  -0.17%  "

# Row 3
$ws.Cells.Item(3, 4).Value = "2.575.59"
$ws.Cells.Item(3, 5).Value = "  +2.35%  "

# Row 4
$ws.Cells.Item(4, 5).Value = "  -0.39%  "

# Row 5
$ws.Cells.Item(5, 4).Value = "'302.97"
$ws.Cells.Item(5, 4).Style = "Normal"
$ws.Cells.Item(5, 5).Value = "  +2.17%  "

# Row 6
$ws.Cells.Item(6, 4).Value = "'97.50"
$ws.Cells.Item(6, 4).Style = "Normal"
$ws.Cells.Item(6, 5).Value = "  +3.48%  "

# Row 7
$ws.Cells.Item(7, 4).Value = "'0.577"
$ws.Cells.Item(7, 4).Style = "Normal"
$ws.Cells.Item(7, 5).Value = "  +0.62%  "

# Row 8
$ws.Cells.Item(8, 5).Value = "  -0.19%  "

# Row 9
$ws.Cells.Item(9, 4).Value = "'0.551"
$ws.Cells.Item(9, 4).Style = "Normal"
$ws.Cells.Item(9, 5).Value = "  +0.35%  "

# Row 10
$ws.Cells.Item(10, 4).Value = "'36.48"
$ws.Cells.Item(10, 4).Style = "Normal"
$ws.Cells.Item(10, 5).Value = "  +0.09%  "

# Row 11
$ws.Cells.Item(11, 5).Value = "  +0.96%  "

# Row 12
$ws.Cells.Item(12, 4).Value = "'7.66"
$ws.Cells.Item(12, 4).Style = "Normal"
$ws.Cells.Item(12, 5).Value = "  -0.01%  "

# Row 13
$ws.Cells.Item(13, 5).Value = "  +6.73%  "

# Row 14
$ws.Cells.Item(14, 4).Value = "2.605.70"
$ws.Cells.Item(14, 5).Value = "  +3.34%  "

# Row 15
$ws.Cells.Item(15, 4).Value = "'0.888"
$ws.Cells.Item(15, 4).Style = "Normal"
$ws.Cells.Item(15, 5).Value = "  +2.68%  "

# Row 16
$ws.Cells.Item(16, 4).Value = "'14.41"
$ws.Cells.Item(16, 4).Style = "Normal"
$ws.Cells.Item(16, 5).Value = "  +2.43%  "

# Row 17
$ws.Cells.Item(17, 4).Value = "43.014.44"
$ws.Cells.Item(17, 5).Value = "  -0.74%  "

# Row 18
$ws.Cells.Item(18, 4).Value = "'12.97"
$ws.Cells.Item(18, 4).Style = "Normal"
$ws.Cells.Item(18, 5).Value = "  +5.58%  "

# Row 19
$ws.Cells.Item(19, 4).Value = "0.0₃0996"
$ws.Cells.Item(19, 5).Value = "  +3.68%  "

# Row 20
$ws.Cells.Item(20, 4).Value = "'6.65"
$ws.Cells.Item(20, 4).Style = "Normal"
$ws.Cells.Item(20, 5).Value = "  +2.11%  "

# Row 21
$ws.Cells.Item(21, 5).Value = "  +0.01%  "

# Row 22
$ws.Cells.Item(22, 4).Value = "'255.07"
$ws.Cells.Item(22, 4).Style = "Normal"
$ws.Cells.Item(22, 5).Value = "  -1.41%  "

# Row 23
$ws.Cells.Item(23, 4).Value = "'2.97"
$ws.Cells.Item(23, 4).Style = "Normal"
$ws.Cells.Item(23, 5).Value = "  +2.70%  "

# Row 24
$ws.Cells.Item(24, 5).Value = "  -0.29%  "

# Row 25
$ws.Cells.Item(25, 4).Value = "'28.74"
$ws.Cells.Item(25, 4).Style = "Normal"
$ws.Cells.Item(25, 5).Value = "  +0.94%  "

# Row 26
$ws.Cells.Item(26, 5).Value = "  -0.12%  "

# Row 27
$ws.Cells.Item(27, 4).Value = "'10.27"
$ws.Cells.Item(27, 4).Style = "Normal"
$ws.Cells.Item(27, 5).Value = "  +2.80%  "

# Row 28
$ws.Cells.Item(28, 4).Value = "'37.77"
$ws.Cells.Item(28, 4).Style = "Normal"
$ws.Cells.Item(28, 5).Value = "  +2.62%  "

# Row 29
$ws.Cells.Item(29, 5).Value = "  -5.03%  "

# Row 30
$ws.Cells.Item(30, 4).Value = "'6.07"
$ws.Cells.Item(30, 4).Style = "Normal"
$ws.Cells.Item(30, 5).Value = "  +0.78%  "

# Row 31
$ws.Cells.Item(31, 4).Value = "'155.58"
$ws.Cells.Item(31, 4).Style = "Normal"
$ws.Cells.Item(31, 5).Value = "  +3.27%  "

# Row 32
$ws.Cells.Item(32, 4).Value = "'2.19"
$ws.Cells.Item(32, 4).Style = "Normal"
$ws.Cells.Item(32, 5).Value = "  +2.13%  "

# Row 33
$ws.Cells.Item(33, 2).Value = "LidoDAOToken"
$ws.Cells.Item(33, 3).Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
$ws.Cells.Item(33, 4).Value = "'3.40"
$ws.Cells.Item(33, 4).Style = "Normal"
$ws.Cells.Item(33, 5).Value = "  -1.27%  "

# Row 34
$ws.Cells.Item(34, 2).Value = "WEMIXToken"
$ws.Cells.Item(34, 3).Value = "https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix"
$ws.Cells.Item(34, 4).Value = "'2.75"
$ws.Cells.Item(34, 4).Style = "Normal"
$ws.Cells.Item(34, 5).Value = "  +0.39%  "

# Row 35
$ws.Cells.Item(35, 4).Value = "'0.0812"
$ws.Cells.Item(35, 4).Style = "Normal"
$ws.Cells.Item(35, 5).Value = "  +2.17%  "

# Row 36
$ws.Cells.Item(36, 4).Value = "'18.39"
$ws.Cells.Item(36, 4).Style = "Normal"
$ws.Cells.Item(36, 5).Value = "  +14.23%  "

# Row 37
$ws.Cells.Item(37, 5).Value = "  +0.91%  "

# Row 38
$ws.Cells.Item(38, 5).Value = "  +1.00%  "

# Row 39
$ws.Cells.Item(39, 4).Value = "'23.75"
$ws.Cells.Item(39, 4).Style = "Normal"
$ws.Cells.Item(39, 5).Value = "  -13.14%  "

# Row 40
$ws.Cells.Item(40, 4).Value = "'3.45"
$ws.Cells.Item(40, 4).Style = "Normal"
$ws.Cells.Item(40, 5).Value = "  -0.52%  "

# Row 41
$ws.Cells.Item(41, 4).Value = "'0.0311"
$ws.Cells.Item(41, 4).Style = "Normal"
$ws.Cells.Item(41, 5).Value = "  +0.81%  "

# Row 42
$ws.Cells.Item(42, 4).Value = "'3.88"
$ws.Cells.Item(42, 4).Style = "Normal"
$ws.Cells.Item(42, 5).Value = "  +2.87%  "

# Row 43
$ws.Cells.Item(43, 5).Value = "  +25.97%  "

# Row 44
$ws.Cells.Item(44, 4).Value = "2.069.00"
$ws.Cells.Item(44, 5).Value = "  +2.90%  "

# Row 45
$ws.Cells.Item(45, 4).Value = "'0.998"
$ws.Cells.Item(45, 4).Style = "Normal"
$ws.Cells.Item(45, 5).Value = "  -0.52%  "

# Row 46
$ws.Cells.Item(46, 4).Value = "'9.28"
$ws.Cells.Item(46, 4).Style = "Normal"
$ws.Cells.Item(46, 5).Value = "  +4.71%  "

# Row 47
$ws.Cells.Item(47, 4).Value = "'85.67"
$ws.Cells.Item(47, 4).Style = "Normal"
$ws.Cells.Item(47, 5).Value = "  +0.47%  "

# Row 48
$ws.Cells.Item(48, 4).Value = "'77.38"
$ws.Cells.Item(48, 4).Style = "Normal"
$ws.Cells.Item(48, 5).Value = "  +14.80%  "

# Row 49
$ws.Cells.Item(49, 4).Value = "'106.66"
$ws.Cells.Item(49, 4).Style = "Normal"
$ws.Cells.Item(49, 5).Value = "  +3.77%  "

# Row 50
$ws.Cells.Item(50, 4).Value = "2.822.51"
$ws.Cells.Item(50, 5).Value = "  +1.84%  "

# Row 51
$ws.Cells.Item(51, 5).Value = "  +2.66%  "
